$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.350.98"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.358.61"
$ws.Range("E3").Value = "  +3.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.13"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.81"
$ws.Range("E6").Value = "  +2.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.526"
$ws.Range("E7").Value = "  -0.78%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  +3.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.19"
$ws.Range("E10").Value = "  -0.15%  "

$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0813"
$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.98"
$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.721.62"
$ws.Range("E15").Value = "  +3.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.74"
$ws.Range("E16").Value = "  +6.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.374.12"
$ws.Range("E17").Value = "  +4.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.812"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.301.51"
$ws.Range("E19").Value = "  +0.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.01"
$ws.Range("E20").Value = "  -4.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0928"
$ws.Range("E21").Value = "  +1.50%  "

$ws.Range("E22").Value = "  +3.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.31"
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "242.02"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("E25").Value = "  +3.22%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.97"
$ws.Range("E28").Value = "  +8.64%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  +6.25%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.55"
$ws.Range("E30").Value = "  -4.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.60"
$ws.Range("E31").Value = "  -0.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.64"
$ws.Range("E32").Value = "  -2.18%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.29"
$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.27"
$ws.Range("E35").Value = "  +1.55%  "

$ws.Range("E36").Value = "  +6.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.12"
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0740"
$ws.Range("E38").Value = "  +0.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  +5.88%  "

$ws.Range("E40").Value = "  +9.05%  "

$ws.Range("E41").Value = "  +2.16%  "

$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("E43").Value = "  +5.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.90"
$ws.Range("E44").Value = "  +2.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0293"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.991.67"
$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.09"
$ws.Range("E47").Value = "  +3.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.44"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "58.65"
$ws.Range("E49").Value = "  +7.00%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.93"
$ws.Range("E50").Value = "  -3.04%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("E51").Value = "  +3.91%  "
